$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.013.04'
$ws.Range("E2").Value = '  -1.83%  '

$ws.Range("D3").Value = '1.555.88'
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9999'
$ws.Range("E5").Value = '  -0.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '286.99'
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3839'
$ws.Range("E7").Value = '  +3.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3247'
$ws.Range("E8").Value = '  -1.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.51'
$ws.Range("E9").Value = '  -11.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.121'
$ws.Range("E10").Value = '  -3.08%  '

$ws.Range("E11").Value = '  -1.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.31'
$ws.Range("E13").Value = '  -6.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.706'
$ws.Range("E14").Value = '  -2.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.792'
$ws.Range("E15").Value = '  -0.97%  '

$ws.Range("D16").Value = '1.548.60'
$ws.Range("E16").Value = '  -1.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001091'
$ws.Range("E17").Value = '  -1.76%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06626'
$ws.Range("E18").Value = '  -1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.16'
$ws.Range("E19").Value = '  -1.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.395'
$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.91'
$ws.Range("E22").Value = '  -3.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.40'
$ws.Range("E23").Value = '  -4.15%  '

$ws.Range("D24").Value = '22.025.78'
$ws.Range("E24").Value = '  -1.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.309'
$ws.Range("E25").Value = '  -1.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.523'
$ws.Range("E26").Value = '  -3.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.20'
$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.86'
$ws.Range("E28").Value = '  -3.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.862'

$ws.Range("D30").Value = '1.727.11'
$ws.Range("E30").Value = '  -0.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.68'
$ws.Range("E31").Value = '  -2.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.094'
$ws.Range("E32").Value = '  +2.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.872'
$ws.Range("E33").Value = '  -2.58%  '

$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.658'
$ws.Range("E34").Value = '  -16.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08153'
$ws.Range("E35").Value = '  -1.68%  '

$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.239'
$ws.Range("E36").Value = '  -5.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06198'
$ws.Range("E37").Value = '  -2.46%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02296'
$ws.Range("E38").Value = '  -5.13%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.243'
$ws.Range("E39").Value = '  -0.81%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2104'
$ws.Range("E40").Value = '  -4.20%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.218'
$ws.Range("E41").Value = '  -5.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.85'
$ws.Range("E42").Value = '  -4.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9995'
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5924'
$ws.Range("E44").Value = '  -3.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.43'
$ws.Range("E45").Value = '  -2.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.722'
$ws.Range("E46").Value = '  -0.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5730'
$ws.Range("E47").Value = '  -4.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.930'
$ws.Range("E48").Value = '  -4.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '119.17'
$ws.Range("E49").Value = '  -4.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.154'
$ws.Range("E50").Value = '  -3.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06875'
$ws.Range("E51").Value = '  -4.16%  '
